# Lattice multiplication exercises: renumber each problem cell.
# Each cell is a <w:p><w:r>...</w:r></w:p> holding 5 <w:t> runs joined by
# <w:br/> manual line breaks: "A x B", "  b    b", "  ----", "d|    |", "d|    |".
# We rebuild each changed cell paragraph via Range.InsertXML so the run/break
# structure (and the xml:space="preserve" flag on the space-padded lines) is
# reproduced exactly, rather than relying on Range.Text which merges runs and
# drops xml:space, and rather than Find/Replace which (in this host) always
# searches from the top of the document regardless of the range it is called on
# -- fatal here since most of the short "d|    |" run texts repeat across cells.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-CellContent($row, $col, $innerXml) {
    $cell = $tbl.Cell($row, $col)
    $rng = $cell.Range
    $xml = "<w:p $wNs>$innerXml</w:p>"
    $rng.InsertXML($xml)
}

Set-CellContent 1 1 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>41 x 44</w:t><w:br/><w:t xml:space="preserve">  4    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>1|    |</w:t></w:r>'
Set-CellContent 1 2 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>82 x 72</w:t><w:br/><w:t xml:space="preserve">  7    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>2|    |</w:t></w:r>'
Set-CellContent 1 3 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>21 x 21</w:t><w:br/><w:t xml:space="preserve">  2    1</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>1|    |</w:t></w:r>'
Set-CellContent 2 1 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>33 x 51</w:t><w:br/><w:t xml:space="preserve">  5    1</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>3|    |</w:t></w:r>'
Set-CellContent 2 2 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>67 x 62</w:t><w:br/><w:t xml:space="preserve">  6    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>7|    |</w:t></w:r>'
Set-CellContent 2 3 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>31 x 53</w:t><w:br/><w:t xml:space="preserve">  5    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>1|    |</w:t></w:r>'
Set-CellContent 3 1 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>93 x 82</w:t><w:br/><w:t xml:space="preserve">  8    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>3|    |</w:t></w:r>'
Set-CellContent 3 2 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>55 x 45</w:t><w:br/><w:t xml:space="preserve">  4    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>5|    |</w:t></w:r>'
Set-CellContent 3 3 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>54 x 93</w:t><w:br/><w:t xml:space="preserve">  9    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>4|    |</w:t></w:r>'
Set-CellContent 4 1 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>64 x 43</w:t><w:br/><w:t xml:space="preserve">  4    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>4|    |</w:t></w:r>'
Set-CellContent 4 2 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>80 x 50</w:t><w:br/><w:t xml:space="preserve">  5    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>0|    |</w:t></w:r>'
Set-CellContent 4 3 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>68 x 54</w:t><w:br/><w:t xml:space="preserve">  5    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>8|    |</w:t></w:r>'
Set-CellContent 5 1 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>11 x 29</w:t><w:br/><w:t xml:space="preserve">  2    9</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>1|    |</w:t></w:r>'
Set-CellContent 5 2 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>69 x 64</w:t><w:br/><w:t xml:space="preserve">  6    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>9|    |</w:t></w:r>'
Set-CellContent 5 3 '<w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>30 x 76</w:t><w:br/><w:t xml:space="preserve">  7    6</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>0|    |</w:t></w:r>'
